# Updated cryptos list with GitHub Actions — refresh Price (D) and Volume(1h) (E)
# columns for each coin row per the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "39.321.63";  E = "  +1.36%  " },
    @{ Row = 3;  D = "2.153.63";   E = "  +3.12%  " },
    @{ Row = 4;  E = "  -0.06%  " },
    @{ Row = 5;  E = "  +0.04%  " },
    @{ Row = 6;  E = "  +0.82%  " },
    @{ Row = 7;  D = "62.68";      E = "  +3.34%  " },
    @{ Row = 8;  D = "0.999";      E = "  -0.06%  " },
    @{ Row = 9;  E = "  +2.00%  " },
    @{ Row = 10; D = "0.0857";     E = "  +2.26%  " },
    @{ Row = 11; E = "  +0.28%  " },
    @{ Row = 12; D = "16.06";      E = "  +7.26%  " },
    @{ Row = 13; D = "2.468.36";   E = "  +3.00%  " },
    @{ Row = 14; D = "22.25";      E = "  +1.62%  " },
    @{ Row = 15; E = "  +2.39%  " },
    @{ Row = 16; D = "5.54";       E = "  +0.99%  " },
    @{ Row = 17; D = "2.148.03";   E = "  +2.97%  " },
    @{ Row = 18; D = "39.411.00";  E = "  +1.91%  " },
    @{ Row = 19; D = "72.10";      E = "  +0.65%  " },
    @{ Row = 20; D = "6.13";       E = "  +1.74%  " },
    @{ Row = 21; D = "0.0₃0854";   E = "  +1.99%  " },
    @{ Row = 22; D = "228.17";     E = "  +0.61%  " },
    @{ Row = 23; E = "  +0.00%  " },
    @{ Row = 24; D = "2.42";       E = "  +1.31%  " },
    @{ Row = 25; D = "2.33";       E = "  -0.20%  " },
    @{ Row = 26; D = "9.78";       E = "  +3.49%  " },
    @{ Row = 27; D = "171.38";     E = "  +0.20%  " },
    @{ Row = 28; E = "  +0.13%  " },
    @{ Row = 29; E = "  +2.37%  " },
    @{ Row = 30; E = "  -2.67%  " },
    @{ Row = 31; D = "2.56";       E = "  +9.54%  " },
    @{ Row = 32; E = "  +0.60%  " },
    @{ Row = 33; D = "4.62";       E = "  +2.47%  " },
    @{ Row = 34; E = "  +2.29%  " },
    @{ Row = 35; D = "7.14";       E = "  +11.00%  " },
    @{ Row = 36; E = "  +1.03%  " },
    @{ Row = 37; E = "  +0.65%  " },
    @{ Row = 38; D = "3.53";       E = "  -0.56%  " },
    @{ Row = 39; D = "0.998";      E = "  -0.07%  " },
    @{ Row = 40; D = "18.20";      E = "  -0.06%  " },
    @{ Row = 41; E = "  +2.41%  " },
    @{ Row = 42; D = "102.80";     E = "  +1.81%  " },
    @{ Row = 43; D = "1.531.81";   E = "  -0.67%  " },
    @{ Row = 44; E = "  +6.29%  " },
    @{ Row = 45; E = "  +6.85%  " },
    @{ Row = 46; E = "  -0.22%  " },
    @{ Row = 47; E = "  +2.07%  " },
    @{ Row = 48; E = "  -0.64%  " },
    @{ Row = 49; E = "  +1.70%  " },
    @{ Row = 50; D = "2.352.89";   E = "  +3.00%  " },
    @{ Row = 51; D = "2.96";       E = "  -0.31%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        # Price column is stored as text in the source data (e.g. "39.312.29",
        # "18.20" with a trailing zero) — force text so Excel doesn't
        # reinterpret/round it as a number.
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    $ws.Range("E$r").Value = $u.E
}
